# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview sheet + per-locale report sheets).
# - Narrow the "Status" column(s) to their new, tighter report width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update status text wherever it currently reads "Ready for handoff".
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Shrink the status columns to the updated report width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
